$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet
$ws.Name = "NorthStar_ClientRepository (9)"

# --- Clear formatting (date style) on Modified Date / Property Modified Date columns so new values store as plain strings ---
foreach ($addr in @("W2", "W3", "W4", "W5", "W6", "W7", "W8", "AI6", "AI7", "AI8")) {
    $ws.Range($addr).ClearFormats()
}

# --- Set updated / new cell values ---
$ws.Range("A2").Value = "NS_June3"
$ws.Range("B2").Value = 8716172871
$ws.Range("C2").Value = "NS_June10"
$ws.Range("O2").Value = 800
$ws.Range("P2").Value = 126
$ws.Range("T2").Value = 801
$ws.Range("U2").Value = "BusinessUnit"
$ws.Range("W2").Value = "2021-06-10T09:34:49"
$ws.Range("A3").Value = "NS_June3"
$ws.Range("B3").Value = 8716172871
$ws.Range("C3").Value = "NS_June10"
$ws.Range("D3").Value = "Business unit June 10"
$ws.Range("O3").Value = 801
$ws.Range("P3").Value = 62
$ws.Range("R3").Value = 800
$ws.Range("T3").Value = "802#804"
$ws.Range("U3").Value = "ServiceLine#BusinessProcess"
$ws.Range("W3").Value = "2021-06-10T09:34:47"
$ws.Range("A4").Value = "NS_June3"
$ws.Range("B4").Value = 8716172871
$ws.Range("C4").Value = "NS_June10"
$ws.Range("D4").Value = "Business unit June 10"
$ws.Range("E4").Value = "Service Line june 10"
$ws.Range("O4").Value = 802
$ws.Range("P4").Value = 25
$ws.Range("R4").Value = 801
$ws.Range("T4").Value = 803
$ws.Range("W4").Value = "2021-06-10T09:34:47"
$ws.Range("A5").Value = "NS_June3"
$ws.Range("B5").Value = 8716172871
$ws.Range("C5").Value = "NS_June10"
$ws.Range("D5").Value = "Business unit June 10"
$ws.Range("E5").Value = "Service Line june 10"
$ws.Range("F5").Value = "Business Process june 10"
$ws.Range("O5").Value = 803
$ws.Range("P5").Value = 862
$ws.Range("R5").Value = 802
$ws.Range("W5").Value = "2021-06-10T09:34:47"
$ws.Range("A6").Value = "NS_June3"
$ws.Range("B6").Value = 8716172871
$ws.Range("C6").Value = "NS_June10"
$ws.Range("D6").Value = "Business unit June 10"
$ws.Range("F6").Value = "Business Process 2 june 10"
$ws.Range("O6").Value = 804
$ws.Range("P6").Value = 863
$ws.Range("Q6").Value = "BusinessProcess"
$ws.Range("R6").Value = 801
$ws.Range("S6").Value = "BusinessUnit"
$ws.Range("T6").Value = 805
$ws.Range("U6").Value = "BusinessKPI"
$ws.Range("W6").Value = "2021-06-10T09:34:47"
$ws.Range("A7").Value = "NS_June3"
$ws.Range("AI7").Value = "2021-06-10T09:34:08"
$ws.Range("B7").Value = 8716172871
$ws.Range("C7").Value = "NS_June10"
$ws.Range("D7").Value = "Business unit June 10"
$ws.Range("F7").Value = "Business Process 2 june 10"
$ws.Range("I7").Value = "Business KPI June 10"
$ws.Range("O7").Value = 805
$ws.Range("P7").Value = 864
$ws.Range("Q7").Value = "BusinessKPI"
$ws.Range("R7").Value = 804
$ws.Range("S7").Value = "BusinessProcess"
$ws.Range("T7").Value = 806
$ws.Range("U7").Value = "ITKPI"
$ws.Range("W7").Value = "2021-06-10T09:34:47"
$ws.Range("A8").Value = "NS_June3"
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = "NA"
$ws.Range("AE8").Value = "NA"
$ws.Range("AF8").Value = "1 : Very-Low"
$ws.Range("AG8").Value = "N/A"
$ws.Range("AI8").Value = "2021-06-10T09:34:34"
$ws.Range("AL8").Value = $false
$ws.Range("B8").Value = 8716172871
$ws.Range("C8").Value = "NS_June10"
$ws.Range("D8").Value = "Business unit June 10"
$ws.Range("F8").Value = "Business Process 2 june 10"
$ws.Range("I8").Value = "Business KPI June 10"
$ws.Range("L8").Value = "IT KPI June 10"
$ws.Range("O8").Value = 806
$ws.Range("P8").Value = 865
$ws.Range("Q8").Value = "ITKPI"
$ws.Range("R8").Value = 805
$ws.Range("S8").Value = "BusinessKPI"
$ws.Range("V8").Value = "nivetha.ramamurthy@ds.dev.accenture.com"
$ws.Range("W8").Value = "2021-06-10T09:34:47"
$ws.Range("X8").Value = $false
$ws.Range("Y8").Value = "Number"
$ws.Range("Z8").Value = 5

# --- Clear cells that no longer hold data ---
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("AA6").ClearContents()
$ws.Range("AB6").ClearContents()
$ws.Range("AE6").ClearContents()
$ws.Range("AF6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AI6").ClearContents()
$ws.Range("AL6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("Y6").ClearContents()
$ws.Range("Z6").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("L7").ClearContents()
